$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for the 2022-Q4
#    entry and push the previous rows (2022-Q3, 2021-Q4) down by one,
#    renumbering the leading index column as we go.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Rows.Item(2).Insert()

# Re-establish formatting on the newly inserted row by cloning it from the
# row directly below (still carrying the original look of a data row).
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.17

# Fix up the index column for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet by duplicating the existing
#    "2022-Q3" sheet (this carries over all formatting faithfully) and
#    placing the copy immediately before it, then overwrite its data.
# ---------------------------------------------------------------------------
$wsQ3Old = $wb.Worksheets.Item(2)   # "2022-Q3" - also used as a formatting-donor sheet
$wsQ3Old.Copy($wsQ3Old, $null)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# The donor sheet has 3 data rows (2022-Q3 has three funds); 2022-Q4 only
# needs two, so drop the extra one.
$wsQ4.Rows.Item(4).Delete()

# Row 2 - 西部利得量化成长混合A
# (format-donor cells are pulled from row 3 of the old sheet so the
# engine doesn't treat same-row-number copy/paste as a no-op)
$wsQ4.Range("A2").Value = 0

$wsQ4.Range("B2").Value = "'000006"
$wsQ3Old.Range("B3").Copy()
$wsQ4.Range("B2").PasteSpecial(-4122)

$wsQ4.Range("C2").Value = "'西部利得量化成长混合A"
$wsQ3Old.Range("C3").Copy()
$wsQ4.Range("C2").PasteSpecial(-4122)

$wsQ4.Range("D2").Value = "'14.73"
$wsQ3Old.Range("D3").Copy()
$wsQ4.Range("D2").PasteSpecial(-4122)

$wsQ4.Range("E2").Value = "'83.72"
$wsQ3Old.Range("E3").Copy()
$wsQ4.Range("E2").PasteSpecial(-4122)

$wsQ4.Range("F2").Value = "'1.01"
$wsQ3Old.Range("F3").Copy()
$wsQ4.Range("F2").PasteSpecial(-4122)

$wsQ4.Range("G2").Value = "'0.1488"
$wsQ3Old.Range("G3").Copy()
$wsQ4.Range("G2").PasteSpecial(-4122)

$wsQ4.Range("H2").Value = 4

# Row 3 - 西部利得量化成长混合C
# (format-donor cells pulled from row 2 of the old sheet, for the same
# same-row-number reason as above)
$wsQ4.Range("A3").Value = 1

$wsQ4.Range("B3").Value = "'011228"
$wsQ3Old.Range("B2").Copy()
$wsQ4.Range("B3").PasteSpecial(-4122)

$wsQ4.Range("C3").Value = "'西部利得量化成长混合C"
$wsQ3Old.Range("C2").Copy()
$wsQ4.Range("C3").PasteSpecial(-4122)

$wsQ4.Range("D3").Value = "'2.01"
$wsQ3Old.Range("D2").Copy()
$wsQ4.Range("D3").PasteSpecial(-4122)

$wsQ4.Range("E3").Value = "'83.72"
$wsQ3Old.Range("E2").Copy()
$wsQ4.Range("E3").PasteSpecial(-4122)

$wsQ4.Range("F3").Value = "'1.01"
$wsQ3Old.Range("F2").Copy()
$wsQ4.Range("F3").PasteSpecial(-4122)

$wsQ4.Range("G3").Value = "'0.0203"
$wsQ3Old.Range("G2").Copy()
$wsQ4.Range("G3").PasteSpecial(-4122)

$wsQ4.Range("H3").Value = 4

Write-Output "done"
